$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the remaining columns (J3:N3) ---
$ws.Range("J3").Value = 0.420731273459036
$ws.Range("K3").Value = 0.000274601843117833
$ws.Range("L3").Value = 0.000115763326662435
$ws.Range("M3").Value = 56290.9380028253
$ws.Range("N3").Formula = "=SQRT(M3)"

# --- Row 4: brand new data row ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 400
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.05
$ws.Range("I4").Value = 128
$ws.Range("J4").Value = 0.0176456439379044
$ws.Range("K4").Value = 0.00011056729969491
$ws.Range("L4").Value = 0.0000635549375038965
$ws.Range("M4").Value = 79247.7238495054
$ws.Range("N4").Formula = "=SQRT(M4)"

# L4 uses the scientific-notation number format (same as K2/L2)
$ws.Range("L4").NumberFormat = "0.00E+00"

# Update the active selection to match the saved view state
$ws.Range("N6").Select()
